$wb = $excel.ActiveWorkbook

# Add the new worksheet "verifyCanOpenItemDetailsPage" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "verifyCanOpenItemDetailsPage"

# Populate column A (item names) first, then column C (item urls), matching
# the shared-string insertion order used when the workbook was authored.
$newSheet.Range("A1").Value = "itemName"
$newSheet.Range("A2").Value = "Sauce Labs Backpack"
$newSheet.Range("A3").Value = "Test.allTheThings() T-Shirt (Red)"
$newSheet.Range("A4").Value = "Sauce Labs Bolt T-Shirt"

$newSheet.Range("C1").Value = "itemUrl"
$newSheet.Range("C2").Value = "https://www.saucedemo.com/inventory-item.html?id=4"
$newSheet.Range("C3").Value = "https://www.saucedemo.com/inventory-item.html?id=3"
$newSheet.Range("C4").Value = "https://www.saucedemo.com/inventory-item.html?id=1"

# Match page setup used by the other sheets
$newSheet.PageSetup.Orientation = 1

# Make the new sheet the active tab (selection on A2, as in the target)
$newSheet.Activate()
$newSheet.Range("A2").Select()
